# Generate Report for Archive
#
# 1. Update the "Ready for handoff" status text to "In Translation" on every
#    cell that currently holds it (this lets the shared-string table collapse
#    the old text and introduce the new one, updating all dependent cells).
# 2. Shrink the "Status" column(s) from ~17.22 chars to ~13.41 chars on the
#    Overview sheet (columns E and F) and on the zh-cn / de-de sheets
#    (column C).

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
